$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update row 3 in place (Fyrflikig jordstjärna record)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 1046478
$ws.Range("B3").Value = 88856
$ws.Range("P3").Value = "Röhälla, 200 m SO. om Brännvinseken., Öl"
$ws.Range("S3").Value = 10
$ws.Range("AC3").Value = "Leg. Ulla-Britt Andersson & Thomas Gunnarsson"
$ws.Range("AI3").Value = "Under enbuske i betad tallskog på sand."
$ws.Range("AW3").Value = "Tommy Knutsson"
$ws.Range("AX3").Value = "Via Tommy Knutsson"
$ws.Range("AY3").Value = "Tommy Knutsson - Import Fynddatabas 2013"

# ---------------------------------------------------------------------------
# 2) Insert a brand-new row at position 11 (pushes the old rows 11-16 down
#    to 12-17) and populate it with the new "Klibbveronika" record.
# ---------------------------------------------------------------------------
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = 55919941
$ws.Range("B11").Value = 104036
$ws.Range("C11").Value = "Ovaliderad"
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 1656
$ws.Range("F11").Value = "Klibbveronika"
$ws.Range("G11").Value = "Veronica triphyllos"
$ws.Range("H11").Value = "L."
$ws.Range("I11").Value = "'10"
$ws.Range("J11").Value = "plantor/tuvor"
$ws.Range("P11").Value = "Röhälla, Öl"
$ws.Range("Q11").Value = 592032.2303689932
$ws.Range("R11").Value = 6286339.845351437
$ws.Range("S11").Value = 50
$ws.Range("T11").Value = "Kalmar"
$ws.Range("U11").Value = "Mörbylånga"
$ws.Range("V11").Value = "Öland"
$ws.Range("W11").Value = "Glömminge"
$ws.Range("X11").Value = "Hö-Mör-5190"

# Dates are stored as plain text in this sheet ("YYYY-MM-DD"); a leading
# apostrophe forces Excel to keep them as text instead of real date serials.
$ws.Range("Y11").Value = "'2015-03-27"
$ws.Range("Z11").Value = "00:00"
$ws.Range("AA11").Value = "'2015-03-27"
$ws.Range("AB11").Value = "00:00"

$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false

$ws.Range("AW11").Value = "Thomas Gunnarsson"
$ws.Range("AX11").Value = "Pav Johnsson"
$ws.Range("AY11").Value = "Floraväkteri Sverige"
